$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 221, shifting existing rows 221-240 down to 222-241
$ws.Rows.Item(221).Insert()

# Populate the newly inserted row 221 with its data
$ws.Cells.Item(221, 1).Value = 4
$ws.Cells.Item(221, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(221, 3).Value = "Los Lagos"
$ws.Cells.Item(221, 4).Value = 44578
$ws.Cells.Item(221, 5).Value = 10
$ws.Cells.Item(221, 6).Value = 100112045
$ws.Cells.Item(221, 7).Value = "Zapallo"
$ws.Cells.Item(221, 8).Value = "Paine"
$ws.Cells.Item(221, 9).Value = "1a nueva(o)"
$ws.Cells.Item(221, 10).Value = 500
$ws.Cells.Item(221, 11).Value = 500
$ws.Cells.Item(221, 12).Value = 500
$ws.Cells.Item(221, 13).Value = 500
$ws.Cells.Item(221, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(221, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(221, 16).Value = 500
$ws.Cells.Item(221, 17).Value = 1
$ws.Cells.Item(221, 18).Value = "Hortaliza"
